# Updated cryptos list values (Price / Volume(1h)) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.145.69"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "1.679.08"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'214.22"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "'22.77"
$ws.Range("E8").Value = "  +6.52%  "
$ws.Range("E9").Value = "  +2.16%  "
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "1.916.41"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "1.679.12"
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("E14").Value = "  +1.97%  "
$ws.Range("E15").Value = "  +3.29%  "
$ws.Range("D16").Value = "'66.59"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").Value = "27.130.21"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "'235.54"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("E19").Value = "  -3.58%  "
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "'4.54"
$ws.Range("E22").Value = "  +1.72%  "
$ws.Range("E23").Value = "  +2.87%  "
$ws.Range("D24").Value = "'2.09"
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").Value = "'146.92"
$ws.Range("D26").Value = "'7.41"
$ws.Range("E26").Value = "  +2.32%  "
$ws.Range("D27").Value = "'16.33"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  +0.74%  "
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").Value = "1.541.63"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("E34").Value = "  +1.92%  "
$ws.Range("E35").Value = "  -3.00%  "
$ws.Range("D36").Value = "'0.606"
$ws.Range("E36").Value = "  +2.79%  "
$ws.Range("E37").Value = "  +1.94%  "
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("E40").Value = "  +2.83%  "
$ws.Range("D41").Value = "'5.80"
$ws.Range("E41").Value = "  +2.96%  "
$ws.Range("D42").Value = "'69.48"
$ws.Range("E42").Value = "  +2.10%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Value = "1.823.49"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "'0.779"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("D47").Value = "'89.79"
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("E48").Value = "  +3.36%  "
$ws.Range("E49").Value = "  +5.87%  "
$ws.Range("D50").Value = "'8.19"
$ws.Range("E50").Value = "  +3.08%  "
$ws.Range("E51").Value = "  +0.12%  "
